$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "άδιεο" -> "άδειο" in the DNS Servers checkpoint text (row 4, column A)
$ws.Range("A4").Value = "Ρύθμιση DNS Servers ώστε το (άδειο) WordPress site να είναι online στη διέυθυνση trian.gr, με χρήση SSL."

# Update the active window view: select A19 as the active cell (the author had
# scrolled down and clicked on row 19 before saving).
$ws.Range("A19").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 10
$aw.ScrollColumn = 1
